$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("S20")

$ws.Cells.Item(1,1).Value = 'Posted  10/1/2019'

$ws.Cells.Item(13,1).Value = 'EE2040'
$ws.Cells.Item(13,2).Value = 'Circuit Analysis I'
$ws.Cells.Item(13,3).Value = 3
$ws.Cells.Item(13,4).Value = 'Rad'
$ws.Cells.Item(13,5).Value = 'TR 1215-130PM'
$ws.Cells.Item(13,6).Value = 'KHC4073'

$ws.Cells.Item(14,1).Value = 'EE2049-01'
$ws.Cells.Item(14,2).Value = 'Electrical Measurements and Circuits Laboratory'
$ws.Cells.Item(14,3).Value = 1
$ws.Cells.Item(14,4).Value = 'Linker'
$ws.Cells.Item(14,5).Value = 'M 600-830PM'
$ws.Cells.Item(14,6).Value = 'ETC249'

$ws.Cells.Item(15,1).Value = 'EE2049-02'
$ws.Cells.Item(15,2).Value = 'Electrical Measurements and Circuits Laboratory'
$ws.Cells.Item(15,3).Value = 1
$ws.Cells.Item(15,4).Value = 'Linker'
$ws.Cells.Item(15,5).Value = 'W 600-830PM'
$ws.Cells.Item(15,6).Value = 'ETC249'

$ws.Cells.Item(16,1).Value = 'EE2049-03'
$ws.Cells.Item(16,2).Value = 'Electrical Measurements and Circuits Laboratory'
$ws.Cells.Item(16,3).Value = 1
$ws.Cells.Item(16,5).Value = 'T 800-1030AM'
$ws.Cells.Item(16,6).Value = 'ETC249'

$ws.Cells.Item(17,1).Value = 'EE2440-01'
$ws.Cells.Item(17,2).Value = 'Digital Engineering'
$ws.Cells.Item(17,3).Value = 2
$ws.Cells.Item(17,4).Value = 'Wang'
$ws.Cells.Item(17,5).Value = 'TR 800-850AM'
$ws.Cells.Item(17,6).Value = 'ETC256'

$ws.Cells.Item(18,1).Value = 'EE2440-02'
$ws.Cells.Item(18,2).Value = 'Digital Engineering'
$ws.Cells.Item(18,3).Value = 1
$ws.Cells.Item(18,4).Value = 'Wang'
$ws.Cells.Item(18,5).Value = 'TR 905-1020AM'
$ws.Cells.Item(18,6).Value = 'ETC256'

$ws.Cells.Item(19,1).Value = 'EE2440-03'
$ws.Cells.Item(19,2).Value = 'Digital Engineering'
$ws.Cells.Item(19,3).Value = 2
$ws.Cells.Item(19,4).Value = 'Ghaforyfard'
$ws.Cells.Item(19,5).Value = 'WF140-230PM'
$ws.Cells.Item(19,6).Value = 'ETC256'

$ws.Cells.Item(20,1).Value = 'EE2440-04'
$ws.Cells.Item(20,2).Value = 'Digital Engineering'
$ws.Cells.Item(20,3).Value = 1
$ws.Cells.Item(20,4).Value = 'Ghaforyfard'
$ws.Cells.Item(20,5).Value = 'WF 245-400PM'
$ws.Cells.Item(20,6).Value = 'ETC256'

$ws.Cells.Item(21,1).Value = 'EE2449-01'
$ws.Cells.Item(21,2).Value = 'Digital Logic Laboratory'
$ws.Cells.Item(21,3).Value = 1
$ws.Cells.Item(21,4).Value = 'Kim'
$ws.Cells.Item(21,5).Value = 'T 600-830PM'
$ws.Cells.Item(21,6).Value = 'ETC248'

$ws.Cells.Item(22,1).Value = 'EE2449-02'
$ws.Cells.Item(22,2).Value = 'Digital Logic Laboratory'
$ws.Cells.Item(22,3).Value = 1
$ws.Cells.Item(22,4).Value = 'Zhao'
$ws.Cells.Item(22,5).Value = 'M 1055AM-125PM'
$ws.Cells.Item(22,6).Value = 'ETC248'

$ws.Cells.Item(23,1).Value = 'EE2449-03'
$ws.Cells.Item(23,2).Value = 'Digital Logic Laboratory'
$ws.Cells.Item(23,3).Value = 1
$ws.Cells.Item(23,4).Value = 'Zhao'
$ws.Cells.Item(23,5).Value = 'W 1055AM-125PM'
$ws.Cells.Item(23,6).Value = 'ETC248'

$ws.Cells.Item(24,1).Value = 'EE2450-01'
$ws.Cells.Item(24,2).Value = 'Embedded System Programming I'
$ws.Cells.Item(24,3).Value = 2
$ws.Cells.Item(24,4).Value = 'Warter-Perez'
$ws.Cells.Item(24,5).Value = 'MW 800-850AM'
$ws.Cells.Item(24,6).Value = 'ETC256'

$ws.Cells.Item(25,1).Value = 'EE2450-02'
$ws.Cells.Item(25,2).Value = 'Embedded System Programming I'
$ws.Cells.Item(25,3).Value = 1
$ws.Cells.Item(25,4).Value = 'Warter-Perez'
$ws.Cells.Item(25,5).Value = 'MW 905-1020AM'
$ws.Cells.Item(25,6).Value = 'ETC256'

$ws.Cells.Item(26,1).Value = 'EE2450-03'
$ws.Cells.Item(26,2).Value = 'Embedded System Programming I'
$ws.Cells.Item(26,3).Value = 2
$ws.Cells.Item(26,4).Value = 'Ghaforyfard'
$ws.Cells.Item(26,5).Value = 'WF 430-520PM'
$ws.Cells.Item(26,6).Value = 'ETC256'

$ws.Cells.Item(27,1).Value = 'EE2450-04'
$ws.Cells.Item(27,2).Value = 'Embedded System Programming I'
$ws.Cells.Item(27,3).Value = 1
$ws.Cells.Item(27,4).Value = 'Ghaforyfard'
$ws.Cells.Item(27,5).Value = 'WF 535-650PM'
$ws.Cells.Item(27,6).Value = 'ETC256'

$ws.Cells.Item(28,1).Value = 'EE3000'
$ws.Cells.Item(28,2).Value = 'Economics for Engineers'
$ws.Cells.Item(28,3).Value = 3
$ws.Cells.Item(28,4).Value = 'Gharibian'
$ws.Cells.Item(28,5).Value = 'F 400-645PM'
$ws.Cells.Item(28,6).Value = 'ETA332'

$ws.Cells.Item(29,1).Value = 'EE3001'
$ws.Cells.Item(29,2).Value = 'Numerical Analysis and Modeling Using MATLAB'
$ws.Cells.Item(29,3).Value = 1
$ws.Cells.Item(29,4).Value = 'Zhang'
$ws.Cells.Item(29,5).Value = 'T 150-420PM'
$ws.Cells.Item(29,6).Value = 'ETC255G'

$ws.Cells.Item(30,1).Value = 'EE3020'
$ws.Cells.Item(30,2).Value = 'Signals and Systems'
$ws.Cells.Item(30,3).Value = 3
$ws.Cells.Item(30,4).Value = 'Rad'
$ws.Cells.Item(30,5).Value = 'TR 1050AM-1205PM'
$ws.Cells.Item(30,6).Value = 'KHC4077'

$ws.Cells.Item(31,1).Value = 'EE3030'
$ws.Cells.Item(31,2).Value = 'Circuit Analysis II'
$ws.Cells.Item(31,3).Value = 3
$ws.Cells.Item(31,4).Value = 'Tabrizi'
$ws.Cells.Item(31,5).Value = 'TR 1215-130PM'
$ws.Cells.Item(31,6).Value = 'SH358B'

$ws.Cells.Item(32,1).Value = 'EE3200'
$ws.Cells.Item(32,2).Value = 'Analog Communication Systems'
$ws.Cells.Item(32,3).Value = 3
$ws.Cells.Item(32,4).Value = 'Tabrizi'
$ws.Cells.Item(32,5).Value = 'TR 925-1040AM'
$ws.Cells.Item(32,6).Value = 'SH358B'

$ws.Cells.Item(33,1).Value = 'EE3300-01'
$ws.Cells.Item(33,2).Value = 'Electric Machines'
$ws.Cells.Item(33,3).Value = 3
$ws.Cells.Item(33,4).Value = 'Rad'
$ws.Cells.Item(33,5).Value = 'TR 925-1040AM'
$ws.Cells.Item(33,6).Value = 'ETA332'

$ws.Cells.Item(34,1).Value = 'EE3300-02'
$ws.Cells.Item(34,2).Value = 'Electric Machines'
$ws.Cells.Item(34,3).Value = 3
$ws.Cells.Item(34,4).Value = 'Zeleke'
$ws.Cells.Item(34,5).Value = 'MW 430-545pm'
$ws.Cells.Item(34,6).Value = 'KHC3097'

$ws.Cells.Item(35,1).Value = 'EE3309-01'
$ws.Cells.Item(35,2).Value = 'Electromagnetic Energy Conversion Laboratory'
$ws.Cells.Item(35,3).Value = 1
$ws.Cells.Item(35,4).Value = 'Fragoso'
$ws.Cells.Item(35,5).Value = 'T 600-830PM'
$ws.Cells.Item(35,6).Value = 'ETA209'

$ws.Cells.Item(36,1).Value = 'EE3309-02'
$ws.Cells.Item(36,2).Value = 'Electromagnetic Energy Conversion Laboratory'
$ws.Cells.Item(36,3).Value = 1
$ws.Cells.Item(36,4).Value = 'Fragoso'
$ws.Cells.Item(36,5).Value = 'M 600-830PM'
$ws.Cells.Item(36,6).Value = 'ETA209'

$ws.Cells.Item(37,1).Value = 'EE3445-01'
$ws.Cells.Item(37,2).Value = 'Computer Organization for CS students'
$ws.Cells.Item(37,3).Value = 3
$ws.Cells.Item(37,4).Value = 'Lim'
$ws.Cells.Item(37,5).Value = 'MW 140-255PM'
$ws.Cells.Item(37,6).Value = 'KHC4075'

$ws.Cells.Item(38,1).Value = 'EE3445-02'
$ws.Cells.Item(38,2).Value = 'Computer Organization for CS students'
$ws.Cells.Item(38,3).Value = 3
$ws.Cells.Item(38,4).Value = 'Lim'
$ws.Cells.Item(38,5).Value = 'MW 1215-130PM'
$ws.Cells.Item(38,6).Value = 'ETA332'

$ws.Cells.Item(39,1).Value = 'EE3450-01'
$ws.Cells.Item(39,2).Value = 'Embedded Systems Programming II'
$ws.Cells.Item(39,3).Value = 2
$ws.Cells.Item(39,4).Value = 'Wang'
$ws.Cells.Item(39,5).Value = 'TR 1050AM-1140AM'
$ws.Cells.Item(39,6).Value = 'ETC256'

$ws.Cells.Item(40,1).Value = 'EE3450-02'
$ws.Cells.Item(40,2).Value = 'Embedded Systems Programming II'
$ws.Cells.Item(40,3).Value = 1
$ws.Cells.Item(40,4).Value = 'Wang'
$ws.Cells.Item(40,5).Value = 'TR 1155AM-110PM'
$ws.Cells.Item(40,6).Value = 'ETC256'

$ws.Cells.Item(41,1).Value = 'EE3450-03'
$ws.Cells.Item(41,2).Value = 'Embedded Systems Programming II'
$ws.Cells.Item(41,3).Value = 2
$ws.Cells.Item(41,4).Value = 'Lin'
$ws.Cells.Item(41,5).Value = 'TR 600-650PM'
$ws.Cells.Item(41,6).Value = 'ETC256'

$ws.Cells.Item(42,1).Value = 'EE3450-04'
$ws.Cells.Item(42,2).Value = 'Embedded Systems Programming II'
$ws.Cells.Item(42,3).Value = 1
$ws.Cells.Item(42,4).Value = 'Lin'
$ws.Cells.Item(42,5).Value = 'TR 705-820PM'
$ws.Cells.Item(42,6).Value = 'ETC256'

$ws.Cells.Item(43,1).Value = 'EE3450-05'
$ws.Cells.Item(43,2).Value = 'Embedded Systems Programming II'
$ws.Cells.Item(43,3).Value = 2
$ws.Cells.Item(43,4).Value = 'Lin'
$ws.Cells.Item(43,5).Value = 'MW 1100-1150AM'
$ws.Cells.Item(43,6).Value = 'ETC256'

$ws.Cells.Item(44,1).Value = 'EE3450-06'
$ws.Cells.Item(44,2).Value = 'Embedded Systems Programming II'
$ws.Cells.Item(44,3).Value = 1
$ws.Cells.Item(44,4).Value = 'Lin'
$ws.Cells.Item(44,5).Value = 'MW 1205-120PM'
$ws.Cells.Item(44,6).Value = 'ETC256'

$ws.Cells.Item(45,1).Value = 'EE3600'
$ws.Cells.Item(45,2).Value = 'Control Systems I'
$ws.Cells.Item(45,3).Value = 3
$ws.Cells.Item(45,4).Value = 'Tabrizi'
$ws.Cells.Item(45,5).Value = 'TR 1050AM-1205PM'
$ws.Cells.Item(45,6).Value = 'SH358B'

$ws.Cells.Item(46,1).Value = 'EE3700'
$ws.Cells.Item(46,2).Value = 'Electronics I'
$ws.Cells.Item(46,3).Value = 3
$ws.Cells.Item(46,4).Value = 'Tabrizi'
$ws.Cells.Item(46,5).Value = 'MW 1215-130PM'
$ws.Cells.Item(46,6).Value = 'SH358B'

$ws.Cells.Item(47,1).Value = 'EE3720'
$ws.Cells.Item(47,2).Value = 'Digital Electronics'
$ws.Cells.Item(47,3).Value = 3
$ws.Cells.Item(47,4).Value = 'Rad'
$ws.Cells.Item(47,5).Value = 'MW 140-255PM'
$ws.Cells.Item(47,6).Value = 'ETA332'

$ws.Cells.Item(48,1).Value = 'EE3810-01'
$ws.Cells.Item(48,2).Value = 'Sensors, Data Acquisition, and Instrumentation with application to Biomedical Engineering'
$ws.Cells.Item(48,3).Value = 2
$ws.Cells.Item(48,4).Value = 'Zhang'
$ws.Cells.Item(48,5).Value = 'TR 500-550PM'
$ws.Cells.Item(48,6).Value = 'KHC4075'

$ws.Cells.Item(49,1).Value = 'EE3810-02'
$ws.Cells.Item(49,2).Value = 'Sensors, Data Acquisition, and Instrumentation with application to Biomedical Engineering'
$ws.Cells.Item(49,3).Value = 1
$ws.Cells.Item(49,4).Value = 'Zhang'
$ws.Cells.Item(49,5).Value = 'F 800-1030AM'
$ws.Cells.Item(49,6).Value = 'ETC251'

$ws.Cells.Item(50,1).Value = 'EE3810-03'
$ws.Cells.Item(50,2).Value = 'Sensors, Data Acquisition, and Instrumentation with application to Biomedical Engineering'
$ws.Cells.Item(50,3).Value = 1
$ws.Cells.Item(50,4).Value = 'Zhang'
$ws.Cells.Item(50,5).Value = 'F 1055AM-125PM'
$ws.Cells.Item(50,6).Value = 'ETC252'

$ws.Cells.Item(51,1).Value = 'EE3810-04'
$ws.Cells.Item(51,2).Value = 'Sensors, Data Acquisition, and Instrumentation with application to Biomedical Engineering'
$ws.Cells.Item(51,3).Value = 1
$ws.Cells.Item(51,4).Value = 'Zhang'
$ws.Cells.Item(51,5).Value = 'F 600-830PM'
$ws.Cells.Item(51,6).Value = 'ETC252'

$ws.Cells.Item(52,1).Value = 'EE4130'
$ws.Cells.Item(52,2).Value = 'Systems Engineering'
$ws.Cells.Item(52,3).Value = 3
$ws.Cells.Item(52,4).Value = 'Harris'
$ws.Cells.Item(52,5).Value = 'MW 430-545PM'
$ws.Cells.Item(52,6).Value = 'ETA129'

$ws.Cells.Item(53,1).Value = 'EE4220'
$ws.Cells.Item(53,2).Value = 'Digital Signal Processing'
$ws.Cells.Item(53,3).Value = 3
$ws.Cells.Item(53,4).Value = 'Mondin'
$ws.Cells.Item(53,5).Value = 'TR 430-545PM'
$ws.Cells.Item(53,6).Value = 'ETC255D'

$ws.Cells.Item(54,1).Value = 'EE4229'
$ws.Cells.Item(54,2).Value = 'Digital Signal Processing Lab'
$ws.Cells.Item(54,3).Value = 1
$ws.Cells.Item(54,4).Value = 'Emrani'
$ws.Cells.Item(54,5).Value = 'W 600-830PM'
$ws.Cells.Item(54,6).Value = 'ETC252'

$ws.Cells.Item(55,1).Value = 'EE4230'
$ws.Cells.Item(55,2).Value = 'Antennas'
$ws.Cells.Item(55,3).Value = 3
$ws.Cells.Item(55,5).Value = 'TR 725-840PM'
$ws.Cells.Item(55,6).Value = 'ETB12'

$ws.Cells.Item(56,1).Value = 'EE4300'
$ws.Cells.Item(56,2).Value = 'Introduction to Power Systems Engineering'
$ws.Cells.Item(56,3).Value = 3
$ws.Cells.Item(56,4).Value = 'Castaneda'
$ws.Cells.Item(56,5).Value = 'F 650-935PM'
$ws.Cells.Item(56,6).Value = 'ETA332'

$ws.Cells.Item(57,1).Value = 'EE4310'
$ws.Cells.Item(57,2).Value = 'Power Systems Analysis'
$ws.Cells.Item(57,3).Value = 3
$ws.Cells.Item(57,4).Value = 'Samaan'
$ws.Cells.Item(57,5).Value = 'F 400-645PM'
$ws.Cells.Item(57,6).Value = 'ETA227'

$ws.Cells.Item(58,1).Value = 'EE4340'
$ws.Cells.Item(58,2).Value = 'Electromagnetic Energy Conversion'
$ws.Cells.Item(58,3).Value = 3
$ws.Cells.Item(58,4).Value = 'Samaan'
$ws.Cells.Item(58,5).Value = 'W 600-845PM'
$ws.Cells.Item(58,6).Value = 'KHB4017'

$ws.Cells.Item(59,1).Value = 'EE4440'
$ws.Cells.Item(59,2).Value = 'Computer Organization'
$ws.Cells.Item(59,3).Value = 3
$ws.Cells.Item(59,4).Value = 'Ghaforyfard'
$ws.Cells.Item(59,5).Value = 'WF 725-840PM'
$ws.Cells.Item(59,6).Value = 'ETC256'

$ws.Cells.Item(60,1).Value = 'EE4450'
$ws.Cells.Item(60,2).Value = 'Embedded Architectures'
$ws.Cells.Item(60,3).Value = 3
$ws.Cells.Item(60,5).Value = 'MW 430-545PM'
$ws.Cells.Item(60,6).Value = 'ETC255D'

$ws.Cells.Item(61,1).Value = 'EE4480-01'
$ws.Cells.Item(61,2).Value = 'Advanced Digital Design'
$ws.Cells.Item(61,3).Value = 2
$ws.Cells.Item(61,4).Value = 'Lin'
$ws.Cells.Item(61,5).Value = 'TR 140-230PM'
$ws.Cells.Item(61,6).Value = 'ETC256'

$ws.Cells.Item(62,1).Value = 'EE4480-02'
$ws.Cells.Item(62,2).Value = 'Advanced Digital Design'
$ws.Cells.Item(62,3).Value = 1
$ws.Cells.Item(62,4).Value = 'Lin'
$ws.Cells.Item(62,5).Value = 'TR 245-400PM'
$ws.Cells.Item(62,6).Value = 'ETC256'

$ws.Cells.Item(63,1).Value = 'EE4610'
$ws.Cells.Item(63,2).Value = 'Digital Control System'
$ws.Cells.Item(63,3).Value = 3
$ws.Cells.Item(63,4).Value = 'Rad'
$ws.Cells.Item(63,5).Value = 'TR 850-1005PM'
$ws.Cells.Item(63,6).Value = 'ETA332'

$ws.Cells.Item(64,1).Value = 'EE4630'
$ws.Cells.Item(64,2).Value = 'Machine Learning Principles and Applications'
$ws.Cells.Item(64,3).Value = 3
$ws.Cells.Item(64,4).Value = 'Mondin'
$ws.Cells.Item(64,5).Value = 'TR 140-255PM'
$ws.Cells.Item(64,6).Value = 'ETC255D'

$ws.Cells.Item(65,1).Value = 'EE4689'
$ws.Cells.Item(65,2).Value = 'Control Systems Laboratory'
$ws.Cells.Item(65,3).Value = 1
$ws.Cells.Item(65,4).Value = 'Fragoso'
$ws.Cells.Item(65,5).Value = 'F 600-830PM'
$ws.Cells.Item(65,6).Value = 'ETC156'

$ws.Cells.Item(66,1).Value = 'EE4820'
$ws.Cells.Item(66,2).Value = 'Biomed Signal Processing'
$ws.Cells.Item(66,3).Value = 3
$ws.Cells.Item(66,4).Value = 'Vincent, P.'
$ws.Cells.Item(66,5).Value = 'TR 600-715PM'
$ws.Cells.Item(66,6).Value = 'ETC255E'

$ws.Cells.Item(67,1).Value = 'EE4962'
$ws.Cells.Item(67,2).Value = 'Senior Design II'
$ws.Cells.Item(67,3).Value = 3
$ws.Cells.Item(67,4).Value = 'Thorburn'
$ws.Cells.Item(67,5).Value = 'F 200-340PM'
$ws.Cells.Item(67,6).Value = 'KHLH2'

$ws.Cells.Item(68,1).Value = 'EE5150'
$ws.Cells.Item(68,2).Value = 'Systems Performance Analysis'
$ws.Cells.Item(68,3).Value = 3
$ws.Cells.Item(68,5).Value = 'MW 725-840PM'
$ws.Cells.Item(68,6).Value = 'KHB2013'

$ws.Cells.Item(69,1).Value = 'EE5160'
$ws.Cells.Item(69,2).Value = 'Systems Architecture'
$ws.Cells.Item(69,3).Value = 3
$ws.Cells.Item(69,4).Value = 'Karimlou'
$ws.Cells.Item(69,5).Value = 'MW 850-1005PM'
$ws.Cells.Item(69,6).Value = 'ETB12'

$ws.Cells.Item(70,1).Value = 'EE5210'
$ws.Cells.Item(70,2).Value = 'Advanced Digital Communication II'
$ws.Cells.Item(70,3).Value = 3
$ws.Cells.Item(70,4).Value = 'Mondin'
$ws.Cells.Item(70,5).Value = 'TR 600-715PM'
$ws.Cells.Item(70,6).Value = 'ETC255D'

$ws.Cells.Item(71,1).Value = 'EE5220'
$ws.Cells.Item(71,2).Value = 'Principles of Signal Compression'
$ws.Cells.Item(71,3).Value = 3
$ws.Cells.Item(71,4).Value = 'Daneshgaran'
$ws.Cells.Item(71,5).Value = 'TR 725-840PM'
$ws.Cells.Item(71,6).Value = 'ETC255D'

$ws.Cells.Item(72,1).Value = 'EE5360'
$ws.Cells.Item(72,2).Value = 'Renwable Energy Sources in Power'
$ws.Cells.Item(72,3).Value = 3
$ws.Cells.Item(72,4).Value = 'Abu-Jaradeh'
$ws.Cells.Item(72,5).Value = 'F 400-645PM'
$ws.Cells.Item(72,6).Value = 'ETA226'

$ws.Cells.Item(73,1).Value = 'EE5370'
$ws.Cells.Item(73,2).Value = 'Faulted Power Systems'
$ws.Cells.Item(73,3).Value = 3
$ws.Cells.Item(73,4).Value = 'Abed'
$ws.Cells.Item(73,5).Value = 'R 600-845PM'
$ws.Cells.Item(73,6).Value = 'ETA129'

$ws.Cells.Item(74,1).Value = 'EE5440'
$ws.Cells.Item(74,2).Value = 'Computer System Architecture '
$ws.Cells.Item(74,3).Value = 3
$ws.Cells.Item(74,4).Value = 'Zhao'
$ws.Cells.Item(74,5).Value = 'TR 430-545PM'
$ws.Cells.Item(74,6).Value = 'ETA129'

$ws.Cells.Item(75,1).Value = 'EE5600'
$ws.Cells.Item(75,2).Value = 'Linear Sys Analysis'
$ws.Cells.Item(75,3).Value = 3
$ws.Cells.Item(75,4).Value = 'Karimlou'
$ws.Cells.Item(75,5).Value = 'MW 600-715PM'
$ws.Cells.Item(75,6).Value = 'KHB2005'

$ws.Activate()
$ws.Range("B27").Select()
